$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 55.62696966666667
$ws.Range("H2").Value = 166.880909
$ws.Range("I2").Value = 0.2670320042914472
$ws.Range("J2").Value = 0.2670320042914472
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 12.840326
$ws.Range("N2").Value = 38.520978
$ws.Range("O2").Value = 0.3393128690704512
$ws.Range("P2").Value = 0.3393128690704511
$ws.Range("Q2").Value = 714.2684249121113
$ws.Range("R2").Value = 6428.415824209002
$ws.Range("S2").Value = 0.09060739550976399
$ws.Range("T2").Value = 0.09060739550976396

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 55.62696966666667
$ws.Range("H3").Value = 166.880909
$ws.Range("I3").Value = 0.2670320042914472
$ws.Range("J3").Value = 0.2670320042914472
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 9.834223333333334
$ws.Range("N3").Value = 29.50267
$ws.Range("O3").Value = 0.2598749077175229
$ws.Range("P3").Value = 0.2598749077175228
$ws.Range("Q3").Value = 547.048043058559
$ws.Range("R3").Value = 4923.43238752703
$ws.Range("S3").Value = 0.06939491747286501
$ws.Range("T3").Value = 0.069394917472865

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 55.62696966666667
$ws.Range("H4").Value = 166.880909
$ws.Range("I4").Value = 0.2670320042914472
$ws.Range("J4").Value = 0.2670320042914472
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 8.654269333333334
$ws.Range("N4").Value = 25.962808
$ws.Range("O4").Value = 0.228693956617749
$ws.Range("P4").Value = 0.2286939566177489
$ws.Range("Q4").Value = 481.4107776924969
$ws.Range("R4").Value = 4332.696999232472
$ws.Range("S4").Value = 0.06106860560497879
$ws.Range("T4").Value = 0.06106860560497877

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 55.62696966666667
$ws.Range("H5").Value = 166.880909
$ws.Range("I5").Value = 0.2670320042914472
$ws.Range("J5").Value = 0.2670320042914472
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 6.513324000000001
$ws.Range("N5").Value = 19.539972
$ws.Range("O5").Value = 0.1721182665942771
$ws.Range("P5").Value = 0.1721182665942771
$ws.Range("Q5").Value = 362.3164765771721
$ws.Range("R5").Value = 3260.848289194548
$ws.Range("S5").Value = 0.04596108570383946
$ws.Range("T5").Value = 0.04596108570383944

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 63.357325
$ws.Range("H6").Value = 190.071975
$ws.Range("I6").Value = 0.3041408435993349
$ws.Range("J6").Value = 0.3041408435993349
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 12.840326
$ws.Range("N6").Value = 38.520978
$ws.Range("O6").Value = 0.3393128690704512
$ws.Range("P6").Value = 0.3393128690704511
$ws.Range("Q6").Value = 813.52870748795
$ws.Range("R6").Value = 7321.758367391551
$ws.Range("S6").Value = 0.1031989022431977
$ws.Range("T6").Value = 0.1031989022431977

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 63.357325
$ws.Range("H7").Value = 190.071975
$ws.Range("I7").Value = 0.3041408435993349
$ws.Range("J7").Value = 0.3041408435993349
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 9.834223333333334
$ws.Range("N7").Value = 29.50267
$ws.Range("O7").Value = 0.2598749077175229
$ws.Range("P7").Value = 0.2598749077175228
$ws.Range("Q7").Value = 623.0700838525834
$ws.Range("R7").Value = 5607.630754673251
$ws.Range("S7").Value = 0.07903857366350671
$ws.Range("T7").Value = 0.0790385736635067

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 63.357325
$ws.Range("H8").Value = 190.071975
$ws.Range("I8").Value = 0.3041408435993349
$ws.Range("J8").Value = 0.3041408435993349
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 8.654269333333334
$ws.Range("N8").Value = 25.962808
$ws.Range("O8").Value = 0.228693956617749
$ws.Range("P8").Value = 0.2286939566177489
$ws.Range("Q8").Value = 548.3113547895334
$ws.Range("R8").Value = 4934.8021931058
$ws.Range("S8").Value = 0.06955517289179187
$ws.Range("T8").Value = 0.06955517289179186

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 63.357325
$ws.Range("H9").Value = 190.071975
$ws.Range("I9").Value = 0.3041408435993349
$ws.Range("J9").Value = 0.3041408435993349
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 6.513324000000001
$ws.Range("N9").Value = 19.539972
$ws.Range("O9").Value = 0.1721182665942771
$ws.Range("P9").Value = 0.1721182665942771
$ws.Range("Q9").Value = 412.6667854983
$ws.Range("R9").Value = 3714.001069484701
$ws.Range("S9").Value = 0.05234819480083867
$ws.Range("T9").Value = 0.05234819480083865

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 52.65915966666667
$ws.Range("H10").Value = 157.977479
$ws.Range("I10").Value = 0.2527853132096735
$ws.Range("J10").Value = 0.2527853132096735
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 12.840326
$ws.Range("N10").Value = 38.520978
$ws.Range("O10").Value = 0.3393128690704512
$ws.Range("P10").Value = 0.3393128690704511
$ws.Range("Q10").Value = 676.1607770060514
$ws.Range("R10").Value = 6085.446993054463
$ws.Range("S10").Value = 0.08577330988404694
$ws.Range("T10").Value = 0.08577330988404691

$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 52.65915966666667
$ws.Range("H11").Value = 157.977479
$ws.Range("I11").Value = 0.2527853132096735
$ws.Range("J11").Value = 0.2527853132096735
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 9.834223333333334
$ws.Range("N11").Value = 29.50267
$ws.Range("O11").Value = 0.2598749077175229
$ws.Range("P11").Value = 0.2598749077175228
$ws.Range("Q11").Value = 517.861936707659
$ws.Range("R11").Value = 4660.757430368931
$ws.Range("S11").Value = 0.065692559942709
$ws.Range("T11").Value = 0.06569255994270899

$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 52.65915966666667
$ws.Range("H12").Value = 157.977479
$ws.Range("I12").Value = 0.2527853132096735
$ws.Range("J12").Value = 0.2527853132096735
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 8.654269333333334
$ws.Range("N12").Value = 25.962808
$ws.Range("O12").Value = 0.228693956617749
$ws.Range("P12").Value = 0.2286939566177489
$ws.Range("Q12").Value = 455.726550622337
$ws.Range("R12").Value = 4101.538955601032
$ws.Range("S12").Value = 0.05781047345277715
$ws.Range("T12").Value = 0.05781047345277713

$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 52.65915966666667
$ws.Range("H13").Value = 157.977479
$ws.Range("I13").Value = 0.2527853132096735
$ws.Range("J13").Value = 0.2527853132096735
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 6.513324000000001
$ws.Range("N13").Value = 19.539972
$ws.Range("O13").Value = 0.1721182665942771
$ws.Range("P13").Value = 0.1721182665942771
$ws.Range("Q13").Value = 342.9861684767321
$ws.Range("R13").Value = 3086.875516290589
$ws.Range("S13").Value = 0.04350896993014042
$ws.Range("T13").Value = 0.0435089699301404

$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 36.672286
$ws.Range("H14").Value = 110.016858
$ws.Range("I14").Value = 0.1760418388995444
$ws.Range("J14").Value = 0.1760418388995444
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 12.840326
$ws.Range("N14").Value = 38.520978
$ws.Range("O14").Value = 0.3393128690704512
$ws.Range("P14").Value = 0.3393128690704511
$ws.Range("Q14").Value = 470.884107405236
$ws.Range("R14").Value = 4237.956966647124
$ws.Range("S14").Value = 0.05973326143344259
$ws.Range("T14").Value = 0.05973326143344257

$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 36.672286
$ws.Range("H15").Value = 110.016858
$ws.Range("I15").Value = 0.1760418388995444
$ws.Range("J15").Value = 0.1760418388995444
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 9.834223333333334
$ws.Range("N15").Value = 29.50267
$ws.Range("O15").Value = 0.2598749077175229
$ws.Range("P15").Value = 0.2598749077175228
$ws.Range("Q15").Value = 360.6434506678734
$ws.Range("R15").Value = 3245.79105601086
$ws.Range("S15").Value = 0.04574885663844214
$ws.Range("T15").Value = 0.04574885663844213

$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 36.672286
$ws.Range("H16").Value = 110.016858
$ws.Range("I16").Value = 0.1760418388995444
$ws.Range("J16").Value = 0.1760418388995444
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 8.654269333333334
$ws.Range("N16").Value = 25.962808
$ws.Range("O16").Value = 0.228693956617749
$ws.Range("P16").Value = 0.2286939566177489
$ws.Range("Q16").Value = 317.3718401130294
$ws.Range("R16").Value = 2856.346561017264
$ws.Range("S16").Value = 0.04025970466820118
$ws.Range("T16").Value = 0.04025970466820116

$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 36.672286
$ws.Range("H17").Value = 110.016858
$ws.Range("I17").Value = 0.1760418388995444
$ws.Range("J17").Value = 0.1760418388995444
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 6.513324000000001
$ws.Range("N17").Value = 19.539972
$ws.Range("O17").Value = 0.1721182665942771
$ws.Range("P17").Value = 0.1721182665942771
$ws.Range("Q17").Value = 238.858480538664
$ws.Range("R17").Value = 2149.726324847976
$ws.Range("S17").Value = 0.03030001615945857
$ws.Range("T17").Value = 0.03030001615945857
